$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04271373187048222
$ws.Range("C2").Value = 0.306821227259698
$ws.Range("D2").Value = 3.537761648806719
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("G2").Value = 4.381533143997669
